# Se obtiene id del cliente a través del DNI.
# Inserts a new "ClienteDNI" column (col F) before the existing "Cliente"
# column on both sheets, fills the DNI value for the existing data rows on
# "Facturados", and removes the duplicated trailing row that was left over
# on "Facturados". Finally makes "Facturados" the active sheet.

$wb = $excel.ActiveWorkbook

# ---- Sheet "Pendientes" (header-only sheet) ----
$wsPend = $wb.Worksheets.Item("Pendientes")
$wsPend.Range("F1").EntireColumn.Insert()
$wsPend.Range("F1").Value = "ClienteDNI"

# ---- Sheet "Facturados" ----
$wsFact = $wb.Worksheets.Item("Facturados")
$wsFact.Range("F1").EntireColumn.Insert()
$wsFact.Range("F1").Value = "ClienteDNI"

# Fill in the DNI for the two remaining data rows.
$wsFact.Range("F2").Value = 12345
$wsFact.Range("F3").Value = 12345

# Row 4 duplicated row 3's data and is removed entirely.
$wsFact.Rows.Item(4).Delete()

# "Facturados" becomes the active sheet/tab.
$wsFact.Activate()
